$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Update the 12-month-period header labels (shift year labels left, add 1401/12 at the end)
$headerRows = @(8,27,36,45,54,63,72,81,90,99,107,115)
foreach ($r in $headerRows) {
    $ws.Range("E$r").Value = "دوازده ماهه منتهی به 1397/12"
    $ws.Range("F$r").Value = "دوازده ماهه منتهی به 1398/12"
    $ws.Range("G$r").Value = "دوازده ماهه منتهی به 1399/12"
    $ws.Range("H$r").Value = "دوازده ماهه منتهی به 1400/12"
    $ws.Range("I$r").Value = "دوازده ماهه منتهی به 1401/12"
}


# Shift data columns E:I left by one year, and set new value for column I
$ws.Range("E10").Value = 7250403
$ws.Range("F10").Value = 8815088
$ws.Range("G10").Value = 16495641
$ws.Range("H10").Value = 44555576
$ws.Range("I10").Value = 36997537
$ws.Range("E11").Value = 139830
$ws.Range("F11").Value = 183383
$ws.Range("G11").Value = 258981
$ws.Range("H11").Value = 355269
$ws.Range("I11").Value = 639750
$ws.Range("E12").Value = 3188170
$ws.Range("F12").Value = 4456588
$ws.Range("G12").Value = 7000668
$ws.Range("H12").Value = 16068124
$ws.Range("I12").Value = 22215729
$ws.Range("E13").Value = 10578403
$ws.Range("F13").Value = 13455059
$ws.Range("G13").Value = 23755290
$ws.Range("H13").Value = 60978969
$ws.Range("I13").Value = 59853016
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("E15").Value = 10578403
$ws.Range("F15").Value = 13455059
$ws.Range("G15").Value = 23755290
$ws.Range("H15").Value = 60978969
$ws.Range("I15").Value = 59853016
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = -1001
$ws.Range("G17").Value = -3814
$ws.Range("H17").Value = -2
$ws.Range("I17").Value = -1
$ws.Range("E18").Value = 10578403
$ws.Range("F18").Value = 13454058
$ws.Range("G18").Value = 23751476
$ws.Range("H18").Value = 60978967
$ws.Range("I18").Value = 59853015
$ws.Range("E19").Value = 352030
$ws.Range("F19").Value = 278578
$ws.Range("G19").Value = 561575
$ws.Range("H19").Value = 492652
$ws.Range("I19").Value = 2672434
$ws.Range("E20").Value = -278578
$ws.Range("F20").Value = -561575
$ws.Range("G20").Value = -492652
$ws.Range("H20").Value = -2672434
$ws.Range("I20").Value = -1834476
$ws.Range("E21").Value = 10651855
$ws.Range("F21").Value = 13171061
$ws.Range("G21").Value = 23820399
$ws.Range("H21").Value = 58799185
$ws.Range("I21").Value = 60690973
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("E23").Value = 10651855
$ws.Range("F23").Value = 13171061
$ws.Range("G23").Value = 23820399
$ws.Range("H23").Value = 58799185
$ws.Range("I23").Value = 60690973
$ws.Range("E29").Value = "-"
$ws.Range("F29").Value = "-"
$ws.Range("G29").Value = "-"
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("E38").Value = "-"
$ws.Range("F38").Value = "-"
$ws.Range("G38").Value = "-"
$ws.Range("H38").Value = 16551
$ws.Range("I38").Value = 5049
$ws.Range("E39").Value = 266950
$ws.Range("F39").Value = 233
$ws.Range("G39").Value = 200432
$ws.Range("H39").Value = 167148
$ws.Range("I39").Value = 161202
$ws.Range("E40").Value = 706963
$ws.Range("F40").Value = 761
$ws.Range("G40").Value = 785283
$ws.Range("H40").Value = 752499
$ws.Range("I40").Value = 689162
$ws.Range("E41").Value = 973913
$ws.Range("F41").Value = 994
$ws.Range("G41").Value = 985715
$ws.Range("H41").Value = 936198
$ws.Range("I41").Value = 855413
$ws.Range("E47").Value = "-"
$ws.Range("F47").Value = "-"
$ws.Range("G47").Value = "-"
$ws.Range("H47").Value = 16551
$ws.Range("I47").Value = 5049
$ws.Range("E48").Value = 266950
$ws.Range("F48").Value = 233
$ws.Range("G48").Value = 200432
$ws.Range("H48").Value = 167148
$ws.Range("I48").Value = 161202
$ws.Range("E49").Value = 706963
$ws.Range("F49").Value = 761
$ws.Range("G49").Value = 785283
$ws.Range("H49").Value = 752499
$ws.Range("I49").Value = 689162
$ws.Range("E50").Value = 973913
$ws.Range("F50").Value = 994
$ws.Range("G50").Value = 985715
$ws.Range("H50").Value = 936198
$ws.Range("I50").Value = 855413
$ws.Range("E56").Value = "-"
$ws.Range("F56").Value = "-"
$ws.Range("G56").Value = "-"
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("E57").Value = 0
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("E59").Value = 0
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("E65").Value = "-"
$ws.Range("F65").Value = "-"
$ws.Range("G65").Value = "-"
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("E66").Value = 0
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("E67").Value = 0
$ws.Range("F67").Value = 0
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("E74").Value = "-"
$ws.Range("F74").Value = "-"
$ws.Range("G74").Value = "-"
$ws.Range("H74").Value = 1507271
$ws.Range("I74").Value = 466179
$ws.Range("E75").Value = 280377
$ws.Range("F75").Value = 584158
$ws.Range("G75").Value = 724539
$ws.Range("H75").Value = 637960
$ws.Range("I75").Value = 804172
$ws.Range("E76").Value = 6970167
$ws.Range("F76").Value = 8230930
$ws.Range("G76").Value = 15771102
$ws.Range("H76").Value = 42410345
$ws.Range("I76").Value = 35727186
$ws.Range("E77").Value = 7250544
$ws.Range("F77").Value = 8815088
$ws.Range("G77").Value = 16495641
$ws.Range("H77").Value = 44555576
$ws.Range("I77").Value = 36997537
$ws.Range("E83").Value = "-"
$ws.Range("F83").Value = "-"
$ws.Range("G83").Value = "-"
$ws.Range("H83").Value = 1507271
$ws.Range("I83").Value = 466179
$ws.Range("E84").Value = 280377
$ws.Range("F84").Value = 584158
$ws.Range("G84").Value = 724539
$ws.Range("H84").Value = 637960
$ws.Range("I84").Value = 804172
$ws.Range("E85").Value = 6970167
$ws.Range("F85").Value = 8230930
$ws.Range("G85").Value = 15771102
$ws.Range("H85").Value = 42410345
$ws.Range("I85").Value = 35727186
$ws.Range("E86").Value = 7250544
$ws.Range("F86").Value = 8815088
$ws.Range("G86").Value = 16495641
$ws.Range("H86").Value = 44555576
$ws.Range("I86").Value = 36997537
$ws.Range("E92").Value = "-"
$ws.Range("F92").Value = "-"
$ws.Range("G92").Value = "-"
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("E94").Value = 0
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("E95").Value = 0
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("E101").Value = "-"
$ws.Range("F101").Value = "-"
$ws.Range("G101").Value = "-"
$ws.Range("H101").Value = 91068274
$ws.Range("I101").Value = 92330957
$ws.Range("E102").Value = 1050298
$ws.Range("F102").Value = 2507115880
$ws.Range("G102").Value = 3614887
$ws.Range("H102").Value = 3816737
$ws.Range("I102").Value = 4988598
$ws.Range("E103").Value = 9859309
$ws.Range("F103").Value = 10815939553
$ws.Range("G103").Value = 20083336
$ws.Range("H103").Value = 56359337
$ws.Range("I103").Value = 51841492
$ws.Range("E109").Value = "-"
$ws.Range("F109").Value = "-"
$ws.Range("G109").Value = "-"
$ws.Range("H109").Value = 91068274
$ws.Range("I109").Value = 92330957
$ws.Range("E110").Value = 1050298
$ws.Range("F110").Value = 2507115880
$ws.Range("G110").Value = 3614887
$ws.Range("H110").Value = 3816737
$ws.Range("I110").Value = 4988598
$ws.Range("E111").Value = 9859309
$ws.Range("F111").Value = 10815939553
$ws.Range("G111").Value = 20083336
$ws.Range("H111").Value = 56359337
$ws.Range("I111").Value = 51841492
$ws.Range("E117").Value = 0
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("E118").Value = 0
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("E119").Value = 0
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("E120").Value = 0
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("E121").Value = 49439
$ws.Range("F121").Value = 92598
$ws.Range("G121").Value = 155661
$ws.Range("H121").Value = 400420
$ws.Range("I121").Value = 601483
$ws.Range("E122").Value = 1372433
$ws.Range("F122").Value = 2292023
$ws.Range("G122").Value = 2423014
$ws.Range("H122").Value = 11042048
$ws.Range("I122").Value = 11296874
$ws.Range("E123").Value = 277098
$ws.Range("F123").Value = 178106
$ws.Range("G123").Value = 154784
$ws.Range("H123").Value = 132050
$ws.Range("I123").Value = 90558
$ws.Range("E124").Value = 309017
$ws.Range("F124").Value = 392384
$ws.Range("G124").Value = 548442
$ws.Range("H124").Value = 880602
$ws.Range("I124").Value = 1359856
$ws.Range("E125").Value = 0
$ws.Range("F125").Value = 0
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("E126").Value = 1180183
$ws.Range("F126").Value = 1501477
$ws.Range("G126").Value = 3718767
$ws.Range("H126").Value = 3613004
$ws.Range("I126").Value = 8866958
$ws.Range("E127").Value = 3188170
$ws.Range("F127").Value = 4456588
$ws.Range("G127").Value = 7000668
$ws.Range("H127").Value = 16068124
$ws.Range("I127").Value = 22215729
